$d = $word.ActiveDocument

# --- Paragraph "Quality peak Brussel Sprouts" ---
# Replace the leading tab character with 16 literal spaces (keep the rest of the
# paragraph text unchanged).
$d.Content.Find.Execute("^tQuality peak Brussel Sprouts", $false, $false, $false, $false, $false, `
    $true, 1, $false, "                Quality peak Brussel Sprouts", 2) | Out-Null

# --- Paragraph "Citric acid ..." ---
# Replace the leading tab character with 16 literal spaces.
$d.Content.Find.Execute("^tCitric acid", $false, $false, $false, $false, $false, `
    $true, 1, $false, "                Citric acid", 2) | Out-Null

# Drop the trailing period after "blanching water".
$d.Content.Find.Execute("blanching water.", $false, $false, $false, $false, $false, `
    $true, 1, $false, "blanching water", 2) | Out-Null

# Append a literal "\n" marker at the end of both paragraphs (matching the
# convention used throughout the rest of the document), just before the
# paragraph mark.
$pQuality = $d.Paragraphs(2)
$endQuality = $d.Range($pQuality.Range.End - 1, $pQuality.Range.End - 1)
$endQuality.InsertAfter("\n")

$pCitric = $d.Paragraphs(3)
$endCitric = $d.Range($pCitric.Range.End - 1, $pCitric.Range.End - 1)
$endCitric.InsertAfter("\n")
